$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ------------------------------------------------------------------
# 1) Opening paragraph: "As a Data Scientist at Mintek, I participate
#    in developing, testing, and deploying machine learning models..."
#    -> wrap "Mintek" in spellcheck proofErr markers and reword
#    "developing, testing, and deploying" as "the development, testing,
#    and deployment of".
# ------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("As a Data Scientist at Mintek")) {
        $targetPara = $p
        break
    }
}

if ($targetPara -ne $null) {
    $xml = '<w:p xmlns:w="' + $wNs + '" w14:paraId="50CD76C6" w14:textId="43F16275" w:rsidR="001A15CF" w:rsidRPr="001A15CF" w:rsidRDefault="001A15CF" w:rsidP="00E65AA2" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr>' +
        '<w:r><w:t xml:space="preserve">As a Data Scientist at </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Mintek</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve">, I </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">participate in </w:t></w:r>' +
        '<w:r><w:t>the development</w:t></w:r>' +
        '<w:r><w:t>, tes</w:t></w:r>' +
        '<w:r><w:t>t</w:t></w:r>' +
        '<w:r><w:t>ing</w:t></w:r>' +
        '<w:r><w:t>, and deploy</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">ment </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">of </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">machine learning models to optimise organisational operations and enhance decision-making, including fine-tuning algorithms for improved performance. I perform A/B testing on machine learning models to compare algorithm performance and feature variations, achieving a 20% improvement in model selection.</w:t></w:r>' +
        '</w:p>'
    [void]$targetPara.Range.InsertXML($xml)
}

# ------------------------------------------------------------------
# 2) Skills bullet: "Machine Learning: Hands-on experience developing
#    and deploying machine learning models using libraries such as
#    scikit-learn, TensorFlow, or PyTorch." -> wrap "PyTorch" in
#    spellcheck proofErr markers (text unchanged).
# ------------------------------------------------------------------
$targetPara2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("Machine Learning: Hands-on experience")) {
        $targetPara2 = $p
        break
    }
}

if ($targetPara2 -ne $null) {
    $xml2 = '<w:p xmlns:w="' + $wNs + '" w14:paraId="128E6270" w14:textId="77777777" w:rsidR="004F3468" w:rsidRPr="004F3468" w:rsidRDefault="004F3468" w:rsidP="004F3468" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
        '<w:r><w:t xml:space="preserve">Machine Learning: Hands-on experience developing and deploying machine learning models using libraries such as scikit-learn, TensorFlow, or </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>PyTorch</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>.</w:t></w:r>' +
        '</w:p>'
    [void]$targetPara2.Range.InsertXML($xml2)
}
